# Update cryptocurrency price/volume data on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> (D value, E value). Empty string means "leave unchanged".
$updates = @{
    2  = @("310.77", "1.78%")
    3  = @("35.63", "-1.93%")
    4  = @("5.108", "0.95%")
    5  = @("0.08224", "4.85%")
    6  = @("2.079", "-10.09%")
    7  = @("", "-0.49%")
    8  = @("", "11.36%")
    9  = @("0.9282", "-0.15%")
    10 = @("0.1076", "10.61%")
    11 = @("0.1917", "3.36%")
    12 = @("0.09308", "4.56%")
    13 = @("0.03626", "-3.82%")
    14 = @("0.09913", "0.20%")
    15 = @("0.001426", "-0.79%")
    16 = @("0.005859", "3.64%")
    17 = @("3.469", "0.04%")
    18 = @("4.125", "-0.72%")
    19 = @("", "0.23%")
    20 = @("0.1309", "-1.13%")
    21 = @("5.098", "-0.78%")
    22 = @("0.2203", "-2.47%")
    23 = @("0.04546", "-1.15%")
    25 = @("", "1.01%")
    26 = @("0.0001250", "")
    27 = @("0.0004445", "-6.20%")
    39 = @("0.01988", "2.55%")
    40 = @("0.04916", "-4.37%")
    41 = @("0.007866", "0.69%")
    42 = @("0.009926", "26.58%")
    43 = @("0.1385", "-0.01%")
    44 = @("0.002115", "-1.76%")
    45 = @("0.01156", "2.66%")
    46 = @("0.00006555", "5.84%")
    47 = @("0.00000000749", "-0.66%")
    48 = @("177.61", "243.53%")
    49 = @("0.001498", "-21.56%")
    50 = @("0.00002098", "-0.66%")
    51 = @("0.0001999", "-0.66%")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne "") {
        $cell = $ws.Range("D$row")
        $cell.Value = "'" + $dVal
        $cell.Style = "Normal"
    }
    if ($eVal -ne "") {
        $cell = $ws.Range("E$row")
        $cell.Value = "'" + $eVal
        $cell.Style = "Normal"
    }
}
